$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 = "Save", using the same style as the other header cells
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New data column H2:H11 = 0 (default/unstyled numeric cells, like column B-G data)
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
